$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1780821917808219
$ws.Range("C2").Value = 0.5958904109589042
$ws.Range("J2").Value = 0.01712328767123288
$ws.Range("O2").Value = 0.003424657534246575
$ws.Range("P2").Value = 0.1267123287671233
$ws.Range("S2").Value = 0.07876712328767123

# Row 3
$ws.Range("C3").Value = 0.01142857142857143
$ws.Range("J3").Value = 0.02857142857142857
$ws.Range("P3").Value = 0.7485714285714286
$ws.Range("S3").Value = 0.2114285714285714

# Row 4
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333

# Row 5
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

# Row 6
$ws.Range("B6").Value = 0.07048458149779736
$ws.Range("F6").Value = 0.08370044052863436
$ws.Range("J6").Value = 0.2731277533039648
$ws.Range("O6").Value = 0.01762114537444934
$ws.Range("Q6").Value = 0.1277533039647577
$ws.Range("R6").Value = 0.05726872246696035
$ws.Range("S6").Value = 0.3700440528634361

# Row 7
$ws.Range("B7").Value = 0.1
$ws.Range("D7").Value = 0.03125
$ws.Range("F7").Value = 0.0375
$ws.Range("J7").Value = 0.1
$ws.Range("O7").Value = 0.01875
$ws.Range("Q7").Value = 0.1625
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.45

# Row 8
$ws.Range("B8").Value = 0.1180952380952381
$ws.Range("D8").Value = 0.01523809523809524
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.07809523809523809
$ws.Range("O8").Value = 0.02476190476190476
$ws.Range("Q8").Value = 0.1714285714285714
$ws.Range("R8").Value = 0.08380952380952381
$ws.Range("S8").Value = 0.4419047619047619

# Row 9
$ws.Range("B9").Value = 0.1020408163265306
$ws.Range("D9").Value = 0.01360544217687075
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.04761904761904762
$ws.Range("O9").Value = 0.01360544217687075
$ws.Range("Q9").Value = 0.1768707482993197
$ws.Range("R9").Value = 0.1156462585034014
$ws.Range("S9").Value = 0.4829931972789115

# Row 10
$ws.Range("B10").Value = 0.1058631921824104
$ws.Range("D10").Value = 0.0252442996742671
$ws.Range("E10").Value = 0.00244299674267101
$ws.Range("F10").Value = 0.07573289902280131
$ws.Range("J10").Value = 0.1042345276872964
$ws.Range("O10").Value = 0.01302931596091205
$ws.Range("Q10").Value = 0.1864820846905537
$ws.Range("R10").Value = 0.07573289902280131
$ws.Range("S10").Value = 0.4112377850162867

# Row 11
$ws.Range("G11").Value = 0.1254901960784314
$ws.Range("J11").Value = 0.06666666666666667
$ws.Range("K11").Value = 0.1803921568627451
$ws.Range("L11").Value = 0.6039215686274509
$ws.Range("S11").Value = 0.02352941176470588

# Row 12
$ws.Range("G12").Value = 0.7341772151898734
$ws.Range("J12").Value = 0.2151898734177215
$ws.Range("K12").Value = 0.006329113924050633
$ws.Range("L12").Value = 0.0189873417721519
$ws.Range("S12").Value = 0.02531645569620253

# Row 13
$ws.Range("G13").Value = 0.6774193548387096
$ws.Range("J13").Value = 0.2580645161290323
$ws.Range("S13").Value = 0.06451612903225806

# Row 15
$ws.Range("F15").Value = 0.0196078431372549
$ws.Range("H15").Value = 0.1617647058823529
$ws.Range("I15").Value = 0.04901960784313725
$ws.Range("J15").Value = 0.392156862745098
$ws.Range("K15").Value = 0.07352941176470588
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("O15").Value = 0.0392156862745098
$ws.Range("S15").Value = 0.25

# Row 16
$ws.Range("F16").Value = 0.01530612244897959
$ws.Range("H16").Value = 0.1836734693877551
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4336734693877551
$ws.Range("K16").Value = 0.07653061224489796
$ws.Range("M16").Value = 0.00510204081632653
$ws.Range("O16").Value = 0.06122448979591837
$ws.Range("S16").Value = 0.1530612244897959

# Row 17
$ws.Range("F17").Value = 0.02544529262086514
$ws.Range("H17").Value = 0.2239185750636132
$ws.Range("I17").Value = 0.07888040712468193
$ws.Range("J17").Value = 0.3867684478371501
$ws.Range("K17").Value = 0.07633587786259542
$ws.Range("M17").Value = 0.01272264631043257
$ws.Range("O17").Value = 0.05089058524173028
$ws.Range("S17").Value = 0.1450381679389313

# Row 18
$ws.Range("F18").Value = 0.02209944751381215
$ws.Range("H18").Value = 0.1602209944751381
$ws.Range("I18").Value = 0.08839779005524862
$ws.Range("J18").Value = 0.4806629834254144
$ws.Range("K18").Value = 0.06077348066298342
$ws.Range("M18").Value = 0.01657458563535912
$ws.Range("O18").Value = 0.06077348066298342
$ws.Range("S18").Value = 0.1104972375690608

# Row 19
$ws.Range("F19").Value = 0.01317715959004392
$ws.Range("H19").Value = 0.2489019033674963
$ws.Range("I19").Value = 0.05563689604685212
$ws.Range("J19").Value = 0.3777452415812592
$ws.Range("K19").Value = 0.09736456808199122
$ws.Range("M19").Value = 0.01464128843338214
$ws.Range("O19").Value = 0.06661786237188873
$ws.Range("S19").Value = 0.1259150805270864
